$d = $word.ActiveDocument

# Update the date heading paragraph (wdReplaceOne=1, wdFindStop=0 to keep the
# replace scoped to the supplied range and avoid touching later duplicate text)
$d.Paragraphs.Item(1).Range.Find.Execute("2023-08-07 Monday", $true, $false, $false, $false, $false, $true, 0, $false, "2023-08-08 Tuesday", 1) | Out-Null

# Update each table cell in row-major order (old -> new computation results).
# Each cell is addressed directly by (row, column) so the fact that a couple of
# source expressions repeat ("39+59=98" appears twice with two different targets)
# cannot cause the wrong cell to be rewritten.
$t = $d.Tables.Item(1)
$pairs = @(
    @(1, 1, "84-3=81", "31+57=88"),
    @(1, 2, "28+28=56", "1+33=34"),
    @(1, 3, "34+16=50", "82-18=64"),
    @(1, 4, "58-25=33", "89+5=94"),
    @(1, 5, "14+59=73", "67-42=25"),
    @(2, 1, "98-81=17", "35+5=40"),
    @(2, 2, "4+45=49", "35+11=46"),
    @(2, 3, "32-1=31", "45-18=27"),
    @(2, 4, "98-32=66", "90-12=78"),
    @(2, 5, "25+60=85", "9+4=13"),
    @(3, 1, "79-77=2", "96-94=2"),
    @(3, 2, "33+20=53", "87-66=21"),
    @(3, 3, "67-9=58", "57+10=67"),
    @(3, 4, "5+33=38", "73+2=75"),
    @(3, 5, "47-14=33", "39+10=49"),
    @(4, 1, "81-2=79", "27-7=20"),
    @(4, 2, "28+20=48", "97-51=46"),
    @(4, 3, "94+2=96", "16+30=46"),
    @(4, 4, "10+26=36", "76-18=58"),
    @(4, 5, "39-34=5", "13-7=6"),
    @(5, 1, "66+15=81", "83-19=64"),
    @(5, 2, "8+70=78", "35-12=23"),
    @(5, 3, "74-4=70", "24-16=8"),
    @(5, 4, "39+26=65", "55-8=47"),
    @(5, 5, "36+10=46", "7+14=21"),
    @(6, 1, "77-40=37", "26+43=69"),
    @(6, 2, "39+59=98", "98-55=43"),
    @(6, 3, "89-8=81", "95-0=95"),
    @(6, 4, "57-49=8", "77-31=46"),
    @(6, 5, "52-6=46", "54-6=48"),
    @(7, 1, "71+22=93", "13-5=8"),
    @(7, 2, "85-3=82", "70-4=66"),
    @(7, 3, "70-20=50", "30+65=95"),
    @(7, 4, "56-14=42", "27+37=64"),
    @(7, 5, "75-31=44", "34-7=27"),
    @(8, 1, "37+3=40", "30+32=62"),
    @(8, 2, "40+6=46", "55-15=40"),
    @(8, 3, "43-23=20", "32+50=82"),
    @(8, 4, "39+59=98", "32+59=91"),
    @(8, 5, "57-50=7", "1+46=47"),
    @(9, 1, "48-26=22", "45-1=44"),
    @(9, 2, "96-33=63", "31+22=53"),
    @(9, 3, "42+39=81", "78-29=49"),
    @(9, 4, "66-52=14", "16-13=3"),
    @(9, 5, "20+7=27", "31-9=22"),
    @(10, 1, "68-34=34", "56+28=84"),
    @(10, 2, "4+91=95", "33+42=75"),
    @(10, 3, "15+44=59", "77+13=90"),
    @(10, 4, "88-27=61", "88-32=56"),
    @(10, 5, "85-6=79", "20+57=77"),
    @(11, 1, "23+52=75", "0+66=66"),
    @(11, 2, "30+43=73", "80-22=58"),
    @(11, 3, "2+49=51", "3+38=41"),
    @(11, 4, "46-23=23", "71-18=53"),
    @(11, 5, "87-15=72", "44-16=28"),
    @(12, 1, "96-35=61", "23+53=76"),
    @(12, 2, "7+79=86", "83-47=36"),
    @(12, 3, "77-2=75", "84-14=70"),
    @(12, 4, "54-8=46", "23+1=24"),
    @(12, 5, "84+11=95", "64-40=24"),
    @(13, 1, "45+4=49", "49-16=33"),
    @(13, 2, "76-54=22", "85-13=72"),
    @(13, 3, "57+6=63", "63-1=62"),
    @(13, 4, "56-50=6", "4+90=94"),
    @(13, 5, "41-23=18", "84-48=36"),
    @(14, 1, "52-1=51", "71+5=76"),
    @(14, 2, "63-59=4", "63+16=79"),
    @(14, 3, "97-16=81", "13-5=8"),
    @(14, 4, "48+39=87", "7+20=27"),
    @(14, 5, "39-28=11", "48-28=20"),
    @(15, 1, "82-33=49", "95-18=77"),
    @(15, 2, "4+59=63", "34-30=4"),
    @(15, 3, "15+15=30", "21+21=42"),
    @(15, 4, "91-86=5", "80-30=50"),
    @(15, 5, "28-20=8", "65+5=70"),
    @(16, 1, "25+26=51", "3+64=67"),
    @(16, 2, "27+62=89", "30+2=32"),
    @(16, 3, "17-17=0", "84+15=99"),
    @(16, 4, "81-55=26", "2+27=29"),
    @(16, 5, "95-25=70", "19+47=66"),
    @(17, 1, "46+5=51", "99-46=53"),
    @(17, 2, "4+58=62", "22-17=5"),
    @(17, 3, "87-85=2", "57-7=50"),
    @(17, 4, "6+8=14", "30-18=12"),
    @(17, 5, "63-32=31", "82+5=87"),
    @(18, 1, "95-80=15", "60-43=17"),
    @(18, 2, "86-26=60", "14+33=47"),
    @(18, 3, "72+23=95", "83-6=77"),
    @(18, 4, "7+41=48", "87-77=10"),
    @(18, 5, "17+32=49", "44+46=90"),
    @(19, 1, "18+24=42", "71-17=54"),
    @(19, 2, "29+2=31", "54-11=43"),
    @(19, 3, "35+14=49", "60-4=56"),
    @(19, 4, "99-35=64", "92-58=34"),
    @(19, 5, "16+71=87", "86-2=84"),
    @(20, 1, "37-33=4", "49+17=66"),
    @(20, 2, "78-61=17", "64+32=96"),
    @(20, 3, "71-22=49", "3+40=43"),
    @(20, 4, "15+73=88", "36-5=31"),
    @(20, 5, "22+41=63", "95-83=12")
)

foreach ($entry in $pairs) {
    $rowIdx = $entry[0]
    $colIdx = $entry[1]
    $oldVal = $entry[2]
    $newVal = $entry[3]
    $cell = $t.Cell($rowIdx, $colIdx)
    # wdFindStop (0) + wdReplaceOne (1): replace only the first match within this
    # cell's own range, never spilling into the rest of the document/table.
    $cell.Range.Find.Execute($oldVal, $true, $false, $false, $false, $false, $true, 0, $false, $newVal, 1) | Out-Null
}

Write-Host "Done"